$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A18").Value = "MAILING Q1"
$ws.Range("B18").Value = 2023
$ws.Range("C18").Value = "Katechezy"
$ws.Range("D18").Value = "INNE"

$ws.Range("B19").Select()
